# "Add Validacion de Fecha"
# Se adiciona validacion de la existencia del campo de fecha
#
# A new ORCID record (0000-0001-9578-9064) is inserted right under the
# header row, pushing the existing rows down by one. A small formatting
# marker (the same underline "flag" style already used at E7) is also
# stamped on F6, and the active selection is moved to D18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capture the existing data rows (2..6) before they get shifted down ---
$colA = @()
$colB = @()
for ($r = 2; $r -le 6; $r++) {
    $colA += ,$ws.Cells.Item($r, 1).Value()
    $colB += ,$ws.Cells.Item($r, 2).Value()
}

# --- Shift rows 2..6 down to 3..7 (write bottom-up so nothing is clobbered) ---
for ($i = 4; $i -ge 0; $i--) {
    $destRow = $i + 3
    $ws.Cells.Item($destRow, 1).Value = $colA[$i]
    $ws.Cells.Item($destRow, 2).Value = $colB[$i]
}

# --- New row 2 holds the newly-registered ORCID id ---
$ws.Cells.Item(2, 1).Value = "xxxxxx"
$ws.Cells.Item(2, 2).Value = "0000-0001-9578-9064"

# --- Stamp the "existence of date field" marker cell at F6 ---
# (matches the pre-existing underline-only style already present at E7)
$ws.Range("F6").Font.Underline = 2

# --- Move the active selection, like the saved workbook shows ---
$ws.Range("D18").Select() | Out-Null
